$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header value label now clarifies it is imported as character in matlab
$ws.Range("B1").Value = "Value (matlab imports as character)"

# Clarify that the practice 'get ready' trigger also marks practice start
$ws.Range("D12").Value = "practice 'get ready'/practice start"

# Clarify that task 'get ready' also marks the block start
$ws.Range("D18").Value = "task 'get ready' (Block start)"

# Add a space after the "S" in several stimulus trigger value codes
$ws.Range("B17").Value = "S 40"
$ws.Range("B18").Value = "S 90"
$ws.Range("B19").Value = "S  1 through S 25"
$ws.Range("B20").Value = "S 80"
$ws.Range("B21").Value = "S 70"
$ws.Range("B22").Value = "S 50"
$ws.Range("B23").Value = "S 30"

# Move the active selection to B2 (as left by the author after editing)
$ws.Range("B2").Select()
